$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value  = 0.1239759412439982
$ws.Range("C3").Value  = 0.02449280966617213
$ws.Range("C4").Value  = 0.09640633152972077
$ws.Range("C5").Value  = 0.0509324021375486
$ws.Range("C6").Value  = 0.05642321045469301
$ws.Range("C7").Value  = 0.03459464430895638
$ws.Range("C8").Value  = 0.1643618244236488
$ws.Range("C9").Value  = 0.06470733700590128
$ws.Range("C10").Value = 0.09134663100090497
$ws.Range("C11").Value = 0.09805506186326893
$ws.Range("C12").Value = 0.07945235310024075
$ws.Range("C13").Value = 0.05254426393812136
$ws.Range("C14").Value = 0.0627071893268247
